$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 678.3
$ws.Range("I8").Value = 97.875
$ws.Range("J8").Value = 3000
$ws.Range("K8").Value = 293.625
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = -154.625
$ws.Range("N8").Value = -9278
$ws.Range("H17").Value = 35955.33
$ws.Range("J17").Value = 35955.33
$ws.Range("L17").Value = 107865.99
$ws.Range("N17").Value = -108201.99
$ws.Range("H28").Value = 434.8095
$ws.Range("I28").Value = 390.29413
$ws.Range("J28").Value = 624
$ws.Range("K28").Value = 390.29413
$ws.Range("L28").Value = 624
$ws.Range("M28").Value = 94.70587
$ws.Range("N28").Value = -1594
$ws.Range("H29").Value = 280.6
$ws.Range("I29").Value = 201
$ws.Range("J29").Value = 400
$ws.Range("K29").Value = 603
$ws.Range("L29").Value = 1200
$ws.Range("M29").Value = -322
$ws.Range("N29").Value = -1762
$ws.Range("H58").Value = 22482.674
$ws.Range("I58").Value = 284.85715
$ws.Range("J58").Value = 26182.309
$ws.Range("K58").Value = 854.5714499999999
$ws.Range("L58").Value = 78546.927
$ws.Range("M58").Value = -704.5714499999999
$ws.Range("N58").Value = -78846.927
$ws.Range("H70").Value = 2305.9285
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 3611.8572
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 10835.5716
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -11375.5716
$ws.Range("H73").Value = 2305.9285
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 3611.8572
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 10835.5716
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -12707.5716
$ws.Range("H82").Value = 5197
$ws.Range("I82").Value = 3795.5
$ws.Range("K82").Value = 11386.5
$ws.Range("M82").Value = -10980.5
$ws.Range("H85").Value = 5197
$ws.Range("I85").Value = 3795.5
$ws.Range("K85").Value = 11386.5
$ws.Range("M85").Value = -9982.5
$ws.Range("H96").Value = 791.4167
$ws.Range("I96").Value = 578
$ws.Range("J96").Value = 943.8570999999999
$ws.Range("K96").Value = 1734
$ws.Range("L96").Value = 2831.5713
$ws.Range("M96").Value = -361
$ws.Range("N96").Value = -5577.5713
$ws.Range("H115").Value = 2162.9167
$ws.Range("I115").Value = 1318.3334
$ws.Range("J115").Value = 2444.4443
$ws.Range("K115").Value = 3955.0002
$ws.Range("L115").Value = 7333.3329
$ws.Range("M115").Value = -2388.0002
$ws.Range("N115").Value = -10467.3329
$ws.Range("H116").Value = 3558.6223
$ws.Range("I116").Value = 3121.138
$ws.Range("J116").Value = 4351.5625
$ws.Range("K116").Value = 3121.138
$ws.Range("L116").Value = 4351.5625
$ws.Range("M116").Value = 320.8620000000001
$ws.Range("N116").Value = -11235.5625
$ws.Range("H132").Value = 6670338.5
$ws.Range("I132").Value = 7146184
$ws.Range("K132").Value = 21438552
$ws.Range("M132").Value = -21436022
$ws.Range("H138").Value = 4332.5884
$ws.Range("J138").Value = 6040.1626
$ws.Range("L138").Value = 18120.4878
$ws.Range("N138").Value = -28400.4878

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 168430.72
$ws.Range("I21").Value = 209799.6
$ws.Range("J21").Value = 65008.5
$ws.Range("K21").Value = 209799.6
$ws.Range("L21").Value = 65008.5
$ws.Range("M21").Value = -209425.6
$ws.Range("N21").Value = -65756.5
$ws.Range("H32").Value = 1206.46
$ws.Range("I32").Value = 1053.7609
$ws.Range("J32").Value = 2962.5
$ws.Range("K32").Value = 1053.7609
$ws.Range("L32").Value = 2962.5
$ws.Range("M32").Value = -766.7609
$ws.Range("N32").Value = -3536.5
$ws.Range("H61").Value = 2866.139
$ws.Range("I61").Value = 1025.4
$ws.Range("J61").Value = 4180.952
$ws.Range("K61").Value = 1025.4
$ws.Range("L61").Value = 4180.952
$ws.Range("M61").Value = -813.4000000000001
$ws.Range("N61").Value = -4604.952
$ws.Range("H136").Value = 2866.139
$ws.Range("I136").Value = 1025.4
$ws.Range("J136").Value = 4180.952
$ws.Range("K136").Value = 3076.2
$ws.Range("L136").Value = 12542.856
$ws.Range("M136").Value = -526.2000000000003
$ws.Range("N136").Value = -17642.856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2774.125
$ws.Range("I31").Value = 1790.762
$ws.Range("K31").Value = 1790.762
$ws.Range("M31").Value = -1495.762
$ws.Range("H34").Value = 2774.125
$ws.Range("I34").Value = 1790.762
$ws.Range("K34").Value = 1790.762
$ws.Range("M34").Value = -1588.762

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 791
$ws.Range("I10").Value = 88.5
$ws.Range("J10").Value = 2196
$ws.Range("K10").Value = 265.5
$ws.Range("L10").Value = 6588
$ws.Range("M10").Value = -126.5
$ws.Range("N10").Value = -6866
$ws.Range("H87").Value = 15057.143
$ws.Range("I87").Value = 10000
$ws.Range("K87").Value = 30000
$ws.Range("M87").Value = -28752
$ws.Range("H90").Value = 15057.143
$ws.Range("I90").Value = 10000
$ws.Range("K90").Value = 90000
$ws.Range("M90").Value = -83760
$ws.Range("H131").Value = 1774.7587
$ws.Range("I131").Value = 4070
$ws.Range("J131").Value = 1176
$ws.Range("K131").Value = 12210
$ws.Range("L131").Value = 3528
$ws.Range("M131").Value = -7170
$ws.Range("N131").Value = -13608

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2219.7856
$ws.Range("I80").Value = 1970
$ws.Range("J80").Value = 2552.8333
$ws.Range("K80").Value = 1970
$ws.Range("L80").Value = 2552.8333
$ws.Range("M80").Value = -972
$ws.Range("N80").Value = -4548.8333
$ws.Range("H83").Value = 2219.7856
$ws.Range("I83").Value = 1970
$ws.Range("J83").Value = 2552.8333
$ws.Range("K83").Value = 9850
$ws.Range("L83").Value = 12764.1665
$ws.Range("M83").Value = -4858
$ws.Range("N83").Value = -22748.1665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1984.7693
$ws.Range("I82").Value = 1300.75
$ws.Range("J82").Value = 2288.7778
$ws.Range("K82").Value = 1300.75
$ws.Range("L82").Value = 2288.7778
$ws.Range("M82").Value = -939.75
$ws.Range("N82").Value = -3010.7778
$ws.Range("H85").Value = 1984.7693
$ws.Range("I85").Value = 1300.75
$ws.Range("J85").Value = 2288.7778
$ws.Range("K85").Value = 1300.75
$ws.Range("L85").Value = 2288.7778
$ws.Range("M85").Value = -52.75
$ws.Range("N85").Value = -4784.7778
$ws.Range("H100").Value = 2189.3333
$ws.Range("I100").Value = 1200
$ws.Range("J100").Value = 2341.5386
$ws.Range("K100").Value = 1200
$ws.Range("L100").Value = 2341.5386
$ws.Range("M100").Value = -659
$ws.Range("N100").Value = -3423.5386
$ws.Range("H136").Value = 1485.65
$ws.Range("I136").Value = 900.9286
$ws.Range("J136").Value = 2850
$ws.Range("K136").Value = 2702.7858
$ws.Range("L136").Value = 8550
$ws.Range("M136").Value = -152.7857999999997
$ws.Range("N136").Value = -13650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4348.6
$ws.Range("I62").Value = 4170.75
$ws.Range("J62").Value = 4467.1665
$ws.Range("K62").Value = 4170.75
$ws.Range("L62").Value = 4467.1665
$ws.Range("M62").Value = -3546.75
$ws.Range("N62").Value = -5715.1665
$ws.Range("H65").Value = 4348.6
$ws.Range("I65").Value = 4170.75
$ws.Range("J65").Value = 4467.1665
$ws.Range("K65").Value = 20853.75
$ws.Range("L65").Value = 22335.8325
$ws.Range("M65").Value = -17733.75
$ws.Range("N65").Value = -28575.8325
$ws.Range("H107").Value = 1241.3334
$ws.Range("I107").Value = 389.6
$ws.Range("J107").Value = 5500
$ws.Range("K107").Value = 1168.8
$ws.Range("L107").Value = 16500
$ws.Range("M107").Value = 751.1999999999998
$ws.Range("N107").Value = -20340
$ws.Range("H122").Value = 209980.77
$ws.Range("I122").Value = 264544.9
$ws.Range("J122").Value = 2637.1
$ws.Range("K122").Value = 793634.7000000001
$ws.Range("L122").Value = 7911.299999999999
$ws.Range("M122").Value = -791184.7000000001
$ws.Range("N122").Value = -12811.3
$ws.Range("H126").Value = 3707253.5
$ws.Range("I126").Value = 2911.2666
$ws.Range("J126").Value = 8337681.5
$ws.Range("K126").Value = 8733.799800000001
$ws.Range("L126").Value = 25013044.5
$ws.Range("M126").Value = -6263.799800000001
$ws.Range("N126").Value = -25017984.5
$ws.Range("H132").Value = 7924.358
$ws.Range("I132").Value = 1585.5555
$ws.Range("J132").Value = 58634.777
$ws.Range("K132").Value = 4756.666499999999
$ws.Range("L132").Value = 175904.331
$ws.Range("M132").Value = -2226.666499999999
$ws.Range("N132").Value = -180964.331
